$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("K2").Value = 63.7
$ws.Range("N2").Value = 54.83846622768671

# Row 3
$ws.Range("D3").Value = 107300
$ws.Range("E3").Value = 60.3
$ws.Range("F3").Value = 6.77
$ws.Range("K3").Value = 58.5
$ws.Range("N3").Value = 54.83846622768671

# Row 4
$ws.Range("K4").Value = 51.7
$ws.Range("N4").Value = 54.83846622768671

# Row 5
$ws.Range("D5").Value = 537000
$ws.Range("E5").Value = 32.3
$ws.Range("F5").Value = 1.32
$ws.Range("K5").Value = 48.9
$ws.Range("N5").Value = 54.83846622768671

# Row 6
$ws.Range("K6").Value = 40.9
$ws.Range("N6").Value = 54.83846622768671

# Row 7
$ws.Range("D7").Value = 64100
$ws.Range("E7").Value = 31.8
$ws.Range("F7").Value = 0.79
$ws.Range("K7").Value = 39.7
$ws.Range("N7").Value = 54.83846622768671
